$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'24.008.05"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.56%  '

$ws.Range("D3").Value = "'1.657.75"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.98%  '

$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.71%  '

$ws.Range("D5").Value = "'310.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.94%  '

$ws.Range("D6").Value = "'1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.44%  '

$ws.Range("D7").Value = "'0.3900"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("D8").Value = "'0.3872"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.18%  '

$ws.Range("D9").Value = "'51.14"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.38%  '

$ws.Range("D10").Value = "'1.368"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.51%  '

$ws.Range("D11").Value = "'1.002"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.77%  '

$ws.Range("D12").Value = "'0.08497"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.13%  '

$ws.Range("D13").Value = "'24.00"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.97%  '

$ws.Range("D14").Value = "'7.213"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.25%  '

$ws.Range("D15").Value = "'8.016"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +7.88%  '

$ws.Range("D16").Value = "'0.00001315"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.59%  '

$ws.Range("D17").Value = "'1.654.54"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.83%  '

$ws.Range("D18").Value = "'94.61"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.59%  '

$ws.Range("D19").Value = "'0.06990"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.18%  '

$ws.Range("D20").Value = "'19.99"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.12%  '

$ws.Range("D21").Value = "'6.996"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.09%  '

$ws.Range("D22").Value = "'1.002"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.26%  '

$ws.Range("D23").Value = "'13.72"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.11%  '

$ws.Range("D24").Value = "'24.010.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.50%  '

$ws.Range("D25").Value = "'2.489"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.51%  '

$ws.Range("D26").Value = "'3.121"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +10.66%  '

$ws.Range("D27").Value = "'22.29"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.08%  '

$ws.Range("D28").Value = "'153.58"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.46%  '

$ws.Range("D29").Value = "'140.66"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.29%  '

$ws.Range("D30").Value = "'5.314"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.76%  '

$ws.Range("D31").Value = "'7.830"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.72%  '

$ws.Range("D32").Value = "'2.493"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.03%  '

$ws.Range("D33").Value = "'1.834.03"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.89%  '

$ws.Range("D34").Value = "'1.049"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +11.10%  '

$ws.Range("D35").Value = "'0.08167"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.78%  '

$ws.Range("D36").Value = "'0.03008"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.68%  '

$ws.Range("D37").Value = "'11.21"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +9.52%  '

$ws.Range("D38").Value = "'6.687"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.20%  '

$ws.Range("D39").Value = "'0.2709"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.08%  '

$ws.Range("D40").Value = "'0.09145"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.10%  '

$ws.Range("D41").Value = "'13.63"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +5.06%  '

$ws.Range("D42").Value = "'0.7585"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.16%  '

$ws.Range("D43").Value = "'1.424"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.14%  '

$ws.Range("D44").Value = "'16.61"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.73%  '

$ws.Range("D45").Value = "'0.7025"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.04%  '

$ws.Range("D46").Value = "'2.505"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.13%  '

$ws.Range("D47").Value = "'4.097"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.64%  '

$ws.Range("D48").Value = "'1.001"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.40%  '

$ws.Range("D49").Value = "'0.08304"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.53%  '

$ws.Range("D50").Value = "'135.23"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.04%  '

$ws.Range("D51").Value = "'1.238"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.38%  '
